$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 478; existing rows 478.. shift down to 479..
$ws.Rows.Item(478).Insert()

# Populate the new row 478 with the weekly Fruta/hortaliza record
$ws.Cells.Item(478, 1).Value = 11
$ws.Cells.Item(478, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(478, 3).Value = "Bíobío"
$ws.Cells.Item(478, 4).Value = 45021
$ws.Cells.Item(478, 5).Value = 8
$ws.Cells.Item(478, 6).Value = "Fruta"
$ws.Cells.Item(478, 7).Value = 100101
$ws.Cells.Item(478, 8).Value = "Berries"
$ws.Cells.Item(478, 9).Value = 100112025
$ws.Cells.Item(478, 10).Value = "Frutilla"
$ws.Cells.Item(478, 11).Value = "Sin especificar"
$ws.Cells.Item(478, 12).Value = "Primera"
$ws.Cells.Item(478, 13).Value = 220
$ws.Cells.Item(478, 14).Value = 7500
$ws.Cells.Item(478, 15).Value = 8000
$ws.Cells.Item(478, 16).Value = 7727
$ws.Cells.Item(478, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(478, 18).Value = "Región del Maule"
$ws.Cells.Item(478, 19).Value = 1104
$ws.Cells.Item(478, 20).Value = 7
